$wb = $excel.ActiveWorkbook

# --- Update "Metadata" sheet timestamp (A2) ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "05 Nov 2025, 12:49 PM"

# --- Update "Industry Analysis" sheet, column F (1 Year) values ---
$wsIndustry = $wb.Worksheets.Item("Industry Analysis")
$wsIndustry.Range("F2").Value = 18.476
$wsIndustry.Range("F3").Value = -7.7404
$wsIndustry.Range("F4").Value = 30.7972
$wsIndustry.Range("F5").Value = -50.2266
$wsIndustry.Range("F6").Value = 61.9649
$wsIndustry.Range("F7").Value = -9.1713
$wsIndustry.Range("F8").Value = -3.556
$wsIndustry.Range("F9").Value = 38.3509
$wsIndustry.Range("F10").Value = -6.2497
$wsIndustry.Range("F11").Value = 52.6723
$wsIndustry.Range("F12").Value = -6.932
$wsIndustry.Range("F13").Value = 17.5662
$wsIndustry.Range("F14").Value = -35.5106
$wsIndustry.Range("F15").Value = 0.6286
$wsIndustry.Range("F16").Value = -3.1514
$wsIndustry.Range("F17").Value = -20.6354
$wsIndustry.Range("F18").Value = -0.0175
$wsIndustry.Range("F19").Value = -26.9255
$wsIndustry.Range("F20").Value = 44.703
$wsIndustry.Range("F21").Value = 10.0506
$wsIndustry.Range("F22").Value = 84.6016
$wsIndustry.Range("F23").Value = -54.4868
$wsIndustry.Range("F24").Value = -12.8122
$wsIndustry.Range("F25").Value = -9.182700000000001
$wsIndustry.Range("F26").Value = 5.9529
$wsIndustry.Range("F27").Value = -33.2998
$wsIndustry.Range("F28").Value = -20.4441
$wsIndustry.Range("F29").Value = -17.1514
$wsIndustry.Range("F30").Value = 24.527
$wsIndustry.Range("F31").Value = 57.6193
$wsIndustry.Range("F32").Value = -1.527
$wsIndustry.Range("F33").Value = -5.2378
$wsIndustry.Range("F34").Value = 27.4054
$wsIndustry.Range("F35").Value = 6.7961
$wsIndustry.Range("F36").Value = -5.6683
$wsIndustry.Range("F37").Value = 1.4178
$wsIndustry.Range("F38").Value = -22.4272
$wsIndustry.Range("F39").Value = 12.3741
$wsIndustry.Range("F40").Value = -5.138
$wsIndustry.Range("F41").Value = -0.1825
$wsIndustry.Range("F42").Value = 23.2483
$wsIndustry.Range("F43").Value = 14.456
$wsIndustry.Range("F44").Value = -11.1739
$wsIndustry.Range("F45").Value = 27.112
$wsIndustry.Range("F46").Value = -5.6252
$wsIndustry.Range("F47").Value = -36.5148
$wsIndustry.Range("F48").Value = -27.8397
$wsIndustry.Range("F49").Value = -25.4424
$wsIndustry.Range("F50").Value = -49.1173
$wsIndustry.Range("F51").Value = -51.065
$wsIndustry.Range("F52").Value = -35.4517
$wsIndustry.Range("F53").Value = -11.9879
$wsIndustry.Range("F54").Value = -3.0992
$wsIndustry.Range("F55").Value = -15.3441
$wsIndustry.Range("F56").Value = -25.937
$wsIndustry.Range("F57").Value = -29.1486
$wsIndustry.Range("F58").Value = -6.4093
$wsIndustry.Range("F59").Value = -23.3046
$wsIndustry.Range("F60").Value = -11.2657
$wsIndustry.Range("F61").Value = -9.777699999999999
$wsIndustry.Range("F62").Value = -16.0561
$wsIndustry.Range("F63").Value = -9.932499999999999
$wsIndustry.Range("F64").Value = 51.8767
$wsIndustry.Range("F65").Value = -43.5191
$wsIndustry.Range("F66").Value = 13.7315
$wsIndustry.Range("F67").Value = 12.6111
$wsIndustry.Range("F68").Value = 31.7532
$wsIndustry.Range("F69").Value = -19.9577
$wsIndustry.Range("F70").Value = -12.9642
$wsIndustry.Range("F71").Value = 13.2432
$wsIndustry.Range("F72").Value = 2.8232
$wsIndustry.Range("F73").Value = -9.179
$wsIndustry.Range("F74").Value = -14.2931
$wsIndustry.Range("F75").Value = 28.3699
$wsIndustry.Range("F76").Value = 45.5868
